$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the affected Price/Volume columns keep their original text formatting
# (these columns store numeric-looking strings like "1.004" or "27.412.98" as
# literal text, e.g. multi-dot "thousand.thousand.hundred" groupings that are
# not valid numbers). Forcing Text number format before assignment prevents
# Excel from auto-converting/normalizing them into real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.412.98'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.05%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.744.00'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.53%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.59'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.49%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4218'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -8.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3581'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.41'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07416'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.10%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.83%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.44'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.58%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.108'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.96%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.189'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.86%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.740.71'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.54%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001067'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.63'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +6.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06074'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -9.56%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.19%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.82%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.29%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -5.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.454.96'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.48'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.73%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.341'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.38'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.30%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.375'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '152.16'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.938.37'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '125.76'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -5.93%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.197'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -5.17%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.677'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.52%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.09117'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.79%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -10.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.67'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02294'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2140'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.87%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.085'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06050'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6387'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.191'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.421'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -5.27%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.910'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.77'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.74%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5857'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.17'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.946'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -5.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06831'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.64%  '

